$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.076.75'
$ws.Range('D3').Value = '1.678.80'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'215.00"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.96%  '
$ws.Range('D9').Value = "'21.26"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.36%  '
$ws.Range('D10').Value = "'0.0622"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').Value = '1.914.81'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').Value = '1.676.45'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('E14').Value = '  +1.08%  '
$ws.Range('D15').Value = "'0.535"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').Value = "'66.07"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('D17').Value = '27.079.17'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').Value = "'236.98"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').Value = "'8.13"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.10%  '
$ws.Range('D20').Value = '0.0₃0742'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('D24').Value = "'2.13"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('D25').Value = "'146.54"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').Value = "'16.31"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.02%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').Value = '1.553.96'
$ws.Range('E32').Value = '  +5.60%  '
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('E35').Value = '  +2.29%  '
$ws.Range('D36').Value = "'0.601"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.80%  '
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').Value = "'0.927"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.33%  '
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = "'68.15"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.05%  '
$ws.Range('D43').Value = "'5.62"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').Value = '1.823.48'
$ws.Range('D46').Value = "'0.785"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'90.70"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.57"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0107'
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('E50').Value = '  +3.00%  '
$ws.Range('D51').Value = "'8.08"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.00%  '
